$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 294.22223
$ws.Range("I9").Value = 413.08334
$ws.Range("K9").Value = 413.08334
$ws.Range("M9").Value = -244.08334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5217.909
$ws.Range("I116").Value = 5216.1
$ws.Range("K116").Value = 5216.1
$ws.Range("M116").Value = -1774.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1628.4324
$ws.Range("I97").Value = 917.2273
$ws.Range("J97").Value = 2671.5334
$ws.Range("K97").Value = 917.2273
$ws.Range("L97").Value = 2671.5334
$ws.Range("M97").Value = -421.2273
$ws.Range("N97").Value = -3663.5334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3717.9656
$ws.Range("I102").Value = 3263.7407
$ws.Range("K102").Value = 3263.7407
$ws.Range("M102").Value = -1641.7407

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 702.8333
$ws.Range("I110").Value = 702.8333
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 702.8333
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1342.1667
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1390.8
$ws.Range("I132").Value = 1220.0435
$ws.Range("J132").Value = 1718.0834
$ws.Range("K132").Value = 3660.1305
$ws.Range("L132").Value = 5154.2502
$ws.Range("M132").Value = -1130.1305
$ws.Range("N132").Value = -10214.2502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1749.9166
$ws.Range("I64").Value = 1792.4286
$ws.Range("J64").Value = 1690.4
$ws.Range("K64").Value = 1792.4286
$ws.Range("L64").Value = 1690.4
$ws.Range("M64").Value = -1567.4286
$ws.Range("N64").Value = -2140.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1749.9166
$ws.Range("I67").Value = 1792.4286
$ws.Range("J67").Value = 1690.4
$ws.Range("K67").Value = 1792.4286
$ws.Range("L67").Value = 1690.4
$ws.Range("M67").Value = -1012.4286
$ws.Range("N67").Value = -3250.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 579.32355
$ws.Range("I94").Value = 522.3889
$ws.Range("J94").Value = 643.375
$ws.Range("K94").Value = 522.3889
$ws.Range("L94").Value = 643.375
$ws.Range("M94").Value = -71.38890000000004
$ws.Range("N94").Value = -1545.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 15904.543
$ws.Range("J107").Value = 1922.875
$ws.Range("L107").Value = 1922.875
$ws.Range("N107").Value = -5762.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2321.8838
$ws.Range("I134").Value = 2121.275
$ws.Range("J134").Value = 4996.6665
$ws.Range("K134").Value = 6363.825000000001
$ws.Range("L134").Value = 14989.9995
$ws.Range("M134").Value = -3828.825000000001
$ws.Range("N134").Value = -20059.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 826.4231
$ws.Range("I16").Value = 890.55554
$ws.Range("K16").Value = 890.55554
$ws.Range("M16").Value = -603.55554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 6668012
$ws.Range("K31").Value = 6668012
$ws.Range("M31").Value = -6667717

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I34").Value = 6668012
$ws.Range("K34").Value = 6668012
$ws.Range("M34").Value = -6667810

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 32222
$ws.Range("J68").Value = 34999.75
$ws.Range("L68").Value = 34999.75
$ws.Range("N68").Value = -36497.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 32222
$ws.Range("J71").Value = 34999.75
$ws.Range("L71").Value = 104999.25
$ws.Range("N71").Value = -112487.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 169500
$ws.Range("I86").Value = 550000
$ws.Range("J86").Value = 42666.668
$ws.Range("K86").Value = 550000
$ws.Range("L86").Value = 42666.668
$ws.Range("M86").Value = -548877
$ws.Range("N86").Value = -44912.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 169500
$ws.Range("I89").Value = 550000
$ws.Range("J89").Value = 42666.668
$ws.Range("K89").Value = 2750000
$ws.Range("L89").Value = 213333.34
$ws.Range("M89").Value = -2744384
$ws.Range("N89").Value = -224565.34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1778.3334
$ws.Range("I105").Value = 1099.1428
$ws.Range("K105").Value = 1099.1428
$ws.Range("M105").Value = 647.8571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 826.4231
$ws.Range("I113").Value = 890.55554
$ws.Range("K113").Value = 890.55554
$ws.Range("M113").Value = 1279.44446

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 272.25
$ws.Range("I26").Value = 20
$ws.Range("J26").Value = 524.5
$ws.Range("K26").Value = 60
$ws.Range("L26").Value = 1573.5
$ws.Range("M26").Value = 228
$ws.Range("N26").Value = -2149.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6975
$ws.Range("J34").Value = 6975
$ws.Range("L34").Value = 20925
$ws.Range("N34").Value = -21093

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2580.2
$ws.Range("J107").Value = 2555.889
$ws.Range("L107").Value = 7667.667
$ws.Range("N107").Value = -11507.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 9999.333000000001
$ws.Range("J113").Value = 9999.333000000001
$ws.Range("L113").Value = 29997.999
$ws.Range("N113").Value = -34337.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 75967.734
$ws.Range("J121").Value = 2901
$ws.Range("L121").Value = 8703
$ws.Range("N121").Value = -11323

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 775.5
$ws.Range("J134").Value = 2500
$ws.Range("L134").Value = 7500
$ws.Range("N134").Value = -17640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 211.91667
$ws.Range("I2").Value = 250.5
$ws.Range("J2").Value = 173.33333
$ws.Range("K2").Value = 250.5
$ws.Range("L2").Value = 173.33333
$ws.Range("M2").Value = -137.5
$ws.Range("N2").Value = -399.33333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3515
$ws.Range("J97").Value = 4927.6665
$ws.Range("L97").Value = 4927.6665
$ws.Range("N97").Value = -5919.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2135.6155
$ws.Range("I102").Value = 1849.0869
$ws.Range("K102").Value = 1849.0869
$ws.Range("M102").Value = -227.0869

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 773.6
$ws.Range("I107").Value = 1180.4
$ws.Range("J107").Value = 366.8
$ws.Range("K107").Value = 1180.4
$ws.Range("L107").Value = 366.8
$ws.Range("M107").Value = 739.5999999999999
$ws.Range("N107").Value = -4206.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 657.0909
$ws.Range("I61").Value = 608.7778
$ws.Range("J61").Value = 874.5
$ws.Range("K61").Value = 608.7778
$ws.Range("L61").Value = 874.5
$ws.Range("M61").Value = -406.7778
$ws.Range("N61").Value = -1278.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3416.4736
$ws.Range("I93").Value = 2972.5
$ws.Range("K93").Value = 2972.5
$ws.Range("M93").Value = -1724.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2210.625
$ws.Range("I100").Value = 2210.625
$ws.Range("K100").Value = 2210.625
$ws.Range("M100").Value = -1669.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 657.0909
$ws.Range("I113").Value = 608.7778
$ws.Range("J113").Value = 874.5
$ws.Range("K113").Value = 608.7778
$ws.Range("L113").Value = 874.5
$ws.Range("M113").Value = 1561.2222
$ws.Range("N113").Value = -5214.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4304.207
$ws.Range("I122").Value = 3417.85
$ws.Range("J122").Value = 6273.8887
$ws.Range("K122").Value = 10253.55
$ws.Range("L122").Value = 18821.6661
$ws.Range("M122").Value = -7803.549999999999
$ws.Range("N122").Value = -23721.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2123
$ws.Range("J107").Value = 2500
$ws.Range("L107").Value = 7500
$ws.Range("N107").Value = -11340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18761.627
$ws.Range("I132").Value = 22231.195
$ws.Range("J132").Value = 918.1429000000001
$ws.Range("K132").Value = 66693.58499999999
$ws.Range("L132").Value = 2754.4287
$ws.Range("M132").Value = -64163.58499999999
$ws.Range("N132").Value = -7814.4287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12610.774
$ws.Range("I136").Value = 14578.604
$ws.Range("J136").Value = 1022.44446
$ws.Range("K136").Value = 43735.812
$ws.Range("L136").Value = 3067.33338
$ws.Range("M136").Value = -41185.812
$ws.Range("N136").Value = -8167.33338
